# Add VAS.AX and NDQ.AX (and IVV.AX) rows to the Australia section of the
# ETF list, mirroring the existing Australia row (row 15) for the
# Country / Asset-Index columns, and select F19 afterwards (matching the
# post-edit cursor position recorded in the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Existing "Australia" / "S&P/ASX 200" values (row 15) reused for the new rows.
$country = $ws.Range("A15").Value2
$asxIndex = $ws.Range("B15").Value2
$spIndex = $ws.Range("B2").Value2

# Row 16: Australia | S&P/ASX 200 | VAS.AX
$ws.Range("A16").Value = $country
$ws.Range("B16").Value = $asxIndex
$ws.Range("C16").Value = "VAS.AX"

# Row 17: Australia | S&P 500 | IVV.AX
$ws.Range("A17").Value = $country
$ws.Range("B17").Value = $spIndex
$ws.Range("C17").Value = "IVV.AX"

# Row 18: Australia | S&P/ASX 200 | NDQ.AX
$ws.Range("A18").Value = $country
$ws.Range("B18").Value = $asxIndex
$ws.Range("C18").Value = "NDQ.AX"

# Match the recorded selection/cursor position after the edit.
$ws.Range("F19").Select()
